# Adding progress events in production namespace and pumpcontrol room (closes #129)
#
# The sheet documents the Socket.IO interface. We need to:
#  1) fix a long-standing typo "orderNumer" -> "orderNumber" in the two
#     existing "state" event rows,
#  2) add a new "progress" event row right under each existing "state"
#     event row in the /orders namespace (one for the "A ordernumber"
#     room, one for the "allOrders" room),
#  3) reword the "allOrders"/"state" comment ("all states" -> "all orders"),
#  4) add a new "pumpControl" room (with a "progress" event) to the
#     /production namespace, right before the closing /connectionState row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the 3 new rows first (bottom-most first so earlier row
#        numbers used below stay stable while we work).
#    Before insistion the sheet has 11 data rows (rows 1-11).
#    - a new row is needed right before the final row (old row 11, the
#      /connectionState row) => insert at row 11
#    - a new row is needed right after old row 4 (allOrders/state)     => insert at row 5
#    - a new row is needed right after old row 3 (A ordernumber/state) => insert at row 4
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(4).Insert()

# Match the row heights used by the surrounding rows (new row 4 sits next
# to the ht=20.25 header-ish rows, the other two new rows sit in the
# ht=20.05 "normal" data-row band).
$ws.Rows.Item(4).RowHeight = 20.25
$ws.Rows.Item(6).RowHeight = 20.05
$ws.Rows.Item(13).RowHeight = 20.05

# --- 2. Fix the orderNumer -> orderNumber typo in the two pre-existing
#        "state" rows (now rows 3 and 5 after the inserts above).
$ws.Range("D3").Value = "{orderNumber:…,fromState:…,toState:…}"
$ws.Range("D5").Value = "{orderNumber:…,fromState:…,toState:…}"

# --- 3. Reword the allOrders/state comment in E5.
$ws.Range("E5").Value = "Statechanges of all orders (on registering the room, events for all orders will be fired (only toState))"

# --- 4. Populate the new "progress" row for the "A ordernumber" room (row 4).
$ws.Range("A4").Value = "/orders"
$ws.Range("B4").Value = "A ordernumber"
$ws.Range("C4").Value = "progress"
$ws.Range("D4").Value = "{orderNumber:…,progress:…}"
$ws.Range("E4").Value = "The progress of a order (on registering the room, a event with the current progress will be fired)"

# --- 5. Populate the new "progress" row for the "allOrders" room (row 6).
$ws.Range("A6").Value = "/orders"
$ws.Range("B6").Value = "allOrders"
$ws.Range("C6").Value = "progress"
$ws.Range("D6").Value = "{orderNumber:…,progress:…}"
$ws.Range("E6").Value = "Progress changes or all orders (on registering the room, events for all orders will be fired)"

# --- 6. Populate the new "pumpControl" room/progress row in /production
#        (row 13, just before the /connectionState row).
$ws.Range("A13").Value = "/production"
$ws.Range("B13").Value = "pumpControl"
$ws.Range("C13").Value = "progress"
$ws.Range("D13").Value = "progress (number 0-100)"
$ws.Range("E13").Value = "The progress from 0 to 100 of the currently relevant order."
